# "Big fix when updating table mfrr"
# Updates the date ranges / intervals / flags on the three sheets of the
# workbook (main, search, week_week) to their corrected values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "main"
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("main")

$wsMain.Range("C5").Value = "2015-01-01"
$wsMain.Range("D5").Value = "2015-01-03"

$wsMain.Range("D10").Value = $false
$wsMain.Range("E10").Value = $false

$wsMain.Range("F11").Value = 60

$wsMain.Range("E12").Value = $true
$wsMain.Range("F12").Value = 60

$wsMain.Range("F13").Value = 60

$wsMain.Range("E11").Select()

# ---------------------------------------------------------------------
# Sheet "search"
# ---------------------------------------------------------------------
$wsSearch = $wb.Worksheets.Item("search")

# C5 text does not change, but it shares the same underlying string as
# main!C5 ("Inital Date" value) -- rewrite it too so the shared-string
# table doesn't keep the stale "2015-05-01" entry around.
$wsSearch.Range("C5").Value = "2015-01-01"
$wsSearch.Range("D5").Value = "2015-09-17"
$wsSearch.Range("F5").Value = 60
$wsSearch.Range("G5").Value = "aFRR_Energy"

$wsSearch.Range("D6").Select()

# ---------------------------------------------------------------------
# Sheet "week_week"
# ---------------------------------------------------------------------
$wsWeek = $wb.Worksheets.Item("week_week")

$wsWeek.Range("C5").Value = "2015-01-01"
$wsWeek.Range("D5").Value = 412

$wsWeek.Range("D10").Value = $false
$wsWeek.Range("E10").Value = $false

$wsWeek.Range("D11").Value = $false
$wsWeek.Range("E11").Value = $false

$wsWeek.Range("D13").Value = $false
$wsWeek.Range("E13").Value = $false

$wsWeek.Range("D7").Select()
